$wb = $excel.ActiveWorkbook

# Sheets "展览" (Exhibitions) and "全部类型" (All types) both contain the
# same event rows; the "想去人数" (# who want to go) counters for the
# first two events each incremented by 1.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 167
    $ws.Range("F3").Value = 117
}
